# Adds a new block of "TNT" rows above the existing "VSEARCH" block in the
# Error_types_table worksheet, shifting the old VSEARCH rows (144-160) down
# to 152-168, and populates the new TNT rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows starting at row 144 (pushes existing rows down).
$ws.Rows("144:151").Insert()

# New TNT data to place into the freshly inserted rows.
$tntData = @(
    @("TNT", "100 Australian species", "12S", "Genus and species correct", 22, 36, 61.11111111111111),
    @("TNT", "100 Australian species", "12S", "Genus correct, species wrong", 12, 36, 33.33333333333333),
    @("TNT", "100 Australian species", "12S", "Genus and species wrong", 2, 36, 5.555555555555555),
    @("TNT", "Lutjanidae", "12S", "Genus and species correct", 14, 16, 87.5),
    @("TNT", "Lutjanidae", "12S", "Genus correct, species wrong", 2, 16, 12.5),
    @("TNT", "Rottnest", "12S", "Genus and species correct", 41, 48, 85.41666666666666),
    @("TNT", "Rottnest", "12S", "Genus correct, species wrong", 6, 48, 12.5),
    @("TNT", "Rottnest", "12S", "Genus and species wrong", 1, 48, 2.083333333333333)
)

$startRow = 144
for ($i = 0; $i -lt $tntData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $tntData[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
    $ws.Cells.Item($r, 7).Value = $rowVals[6]
}
